$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price (D) and Volume(1h) (E) columns are stored as text; force text
# number format before assignment so Excel does not reinterpret numeric-
# looking strings (e.g. "1.00", "42.628.55") as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.628.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.529.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.86"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.11"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.35"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.44"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.917.91"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.576.07"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.96"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.841"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.677.24"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.46"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0954"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.22"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.34"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.47"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.11"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.17"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.31"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.37"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.70"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.08"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.40"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.118"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.15"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.14%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.79"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0302"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.22"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.985.68"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "83.90"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.98"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.774.57"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.87"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.88"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.09%  "

# Rows 35/36: Hedera and LidoDAOToken swap positions, with updated prices
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0793"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.27"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.27%  "
